$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("reaction2")

# Clear the existing row range that held the old values (A1:O1),
# then set the new reduced values.
$ws.Range("A1:O1").Clear()

$ws.Range("A1").Value = 4
$ws.Range("B1").Value = 5
